$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update current_count (B) and change (D) for the rows whose
# current_count increased by 1 relative to previous_count (C).
# Row 15: Department of State          31 -> 32, change 0 -> 1
# Row 21: Executive Office of the President, Management and Administration 6 -> 7, change 0 -> 1
# Row 27: International Development    19 -> 20, change 0 -> 1
# Row 29: National Security Council    19 -> 20, change 0 -> 1

$ws.Range("B15").Value = 32
$ws.Range("D15").Value = 1

$ws.Range("B21").Value = 7
$ws.Range("D21").Value = 1

$ws.Range("B27").Value = 20
$ws.Range("D27").Value = 1

$ws.Range("B29").Value = 20
$ws.Range("D29").Value = 1
